$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column A. This shifts the existing 9 columns
#    of data (A:I) over to (B:J), matching the header/data layout seen in the
#    target sheet.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).Insert()

# ---------------------------------------------------------------------------
# 2. Fill in the new "NA" values first (so the shared string table interns
#    "NA" right after the pre-existing strings), then rewrite the header row
#    with ID/P1..P9 labels instead of the old plain numbers.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "NA"
$ws.Range("E2").Value = "NA"

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "P1"
$ws.Range("C1").Value = "P2"
$ws.Range("D1").Value = "P3"
$ws.Range("E1").Value = "P4"
$ws.Range("F1").Value = "P5"
$ws.Range("G1").Value = "P6"
$ws.Range("H1").Value = "P7"
$ws.Range("I1").Value = "P8"
$ws.Range("J1").Value = "P9"

# ---------------------------------------------------------------------------
# 3. Formatting: data row (row 2) is centered horizontally; header row
#    (row 1) is centered horizontally + vertically and uses an explicit
#    black font color. Build the combined header format on a scratch cell
#    and paste the format in one shot so only a single new style record is
#    generated for the header (instead of one per property assignment).
# ---------------------------------------------------------------------------
$ws.Range("A2:J2").HorizontalAlignment = -4108

$scratch = $ws.Range("Z1")
$scratch.Font.Color = 0
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("A1:J1").PasteSpecial(-4122)
$scratch.Clear()

# ---------------------------------------------------------------------------
# 4. Column widths (best-fit sizing of the new layout). ColumnWidth takes a
#    character-unit width and the engine stores it with a fixed +5/6 offset,
#    so we back-solve for the character width that reproduces each target
#    stored width as closely as possible.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 2.6666666666666665   # -> 3.5
$ws.Columns.Item(2).ColumnWidth = 8.5                  # -> 9.33203125 (approx)
$ws.Columns.Item(3).ColumnWidth = 6                    # -> 6.83203125 (approx)
$ws.Columns.Item(4).ColumnWidth = 3.1666666666666665   # -> 4
$ws.Columns.Item(5).ColumnWidth = 2.6666666666666665   # -> 3.5
$ws.Columns.Item(6).ColumnWidth = 3.1666666666666665   # -> 4
$ws.Columns.Item(7).ColumnWidth = 31.333333333333332   # -> 32.1640625 (approx)
$ws.Columns.Item(8).ColumnWidth = 3.5                  # -> 4.33203125 (approx)
$ws.Columns.Item(9).ColumnWidth = 3.5                  # -> 4.33203125 (approx)

# ---------------------------------------------------------------------------
# 5. Cosmetic sheet-level touches present in the target: active selection
#    moved to B1, and a portrait A4 page setup.
# ---------------------------------------------------------------------------
$ws.Range("B1").Select() | Out-Null
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
